$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Descripcion Bases de Datos" header (old row 12) and its
# description paragraph (old row 13) down to rows 23/24, leaving rows 11-22
# free to hold the 9 new questions (11-19) plus a 3-row gap (20-22), exactly
# matching the target layout.
$ws.Rows("11:21").Insert()

$ws.Range("A11").Value = "¿Existe un cambio en el desempleo en cuanto al ambito racial?"
$ws.Range("A12").Value = "¿La raza es un factor importante en el empleo o desempleo de una persona?"
$ws.Range("A13").Value = "¿Se puede cuantificar la raza en los datos?"
$ws.Range("A14").Value = "¿Hay una tendencia en el nivel de vida de una persona y el desempleo de cada grupo?"
$ws.Range("A15").Value = "Aunque esten menos desempleados, ¿tiene un mejor nivel de vida?"
$ws.Range("A16").Value = "¿Qué tan importante es el tema del desempleo en USA? ¿Es un tema mayor o menor?"
$ws.Range("A17").Value = "¿Las personas están dispuestas a trabajar?"
$ws.Range("A18").Value = "¿Murieron más personas y por eso bajo el desempleo? ¿O quisieron trabajar más?"
$ws.Range("A19").Value = "¿Los migrantes hacen que aumente el desempleo?"

# The long description paragraph now living in row 24 gets a shorter,
# recalculated row height.
$ws.Rows(24).RowHeight = 72

$ws.Range("A19").Select()
